$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.622.45'
$ws.Range('E2').Value = '  -2.77%  '
$ws.Range('D3').Value = '1.984.58'
$ws.Range('E3').Value = '  -3.78%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.638'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.90'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.60%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.79'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.103'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.949'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.52'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('D15').Value = '2.272.14'
$ws.Range('E15').Value = '  -3.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.92%  '
$ws.Range('D17').Value = '1.984.91'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.36%  '
$ws.Range('D19').Value = '35.573.54'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +19.92%  '
$ws.Range('E26').Value = '  -3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('E29').Value = '  -4.70%  '
$ws.Range('E30').Value = '  -2.82%  '
$ws.Range('E31').Value = '  -4.51%  '
$ws.Range('E32').Value = '  -6.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0959'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +16.74%  '
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('E35').Value = '  +9.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.70%  '
$ws.Range('E39').Value = '  +8.00%  '
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('E41').Value = '  +1.68%  '
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0890'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.32%  '
$ws.Range('D48').Value = '1.374.92'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.25'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.84%  '
